$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new profit row for 2025-09-12.
# Force the date to be stored as text (matching the existing rows) rather
# than letting Excel auto-convert the "MM/DD/YYYY" looking string into a
# date serial number, then reset the style so no extra formatting sticks.
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "09/12/2025"
$ws.Range("A26").Style = "Normal"

$ws.Range("B26").Value = 15960.33
